$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.134.44"
$ws.Range("E2").Value = "  -1.00%  "

$ws.Range("D3").Value = "1.671.26"
$ws.Range("E3").Value = "  -1.34%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.72%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.93%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5252"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.51%  "

$ws.Range("E7").Value = "  -0.71%  "

$ws.Range("E8").Value = "  -3.29%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06285"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.74%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.23"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.64%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07528"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.83%  "

$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.446"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.14%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.636.75"
$ws.Range("E13").Value = "  -3.36%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5636"
$ws.Range("D14").Style = "Normal"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008025"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.01%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.46%  "

$ws.Range("D17").Value = "26.196.82"
$ws.Range("E17").Value = "  -0.90%  "

$ws.Range("E18").Value = "  -0.71%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.805"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.63%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "187.92"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.14%  "

$ws.Range("E21").Value = "  -5.42%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.177"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.22%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.004"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.63%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "148.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.32%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1250"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.93%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.596"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.06%  "

$ws.Range("E27").Value = "  +1.10%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06238"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.72%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.357"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.01%  "

$ws.Range("E30").Value = "  -4.09%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.480"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.60%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.436"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.64%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.629"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.22%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9989"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.24%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6045"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.72%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.404"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.31%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.715"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.31%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.115"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.03%  "

$ws.Range("E39").Value = "  -1.80%  "

$ws.Range("D40").Value = "1.076.68"
$ws.Range("E40").Value = "  -3.73%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8665"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.02%  "

$ws.Range("E42").Value = "  -1.10%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.02"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.68%  "

$ws.Range("D44").Value = "1.821.56"
$ws.Range("E44").Value = "  -1.32%  "

$ws.Range("E45").Value = "  +0.95%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.15"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.42%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.000"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.26%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.011"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.34%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05241"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.74%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4255"
$ws.Range("D50").Style = "Normal"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.994"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.85%  "
